$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.175.52'
$ws.Range("E2").Value = '  -3.77%  '
$ws.Range("D3").Value = '2.239.86'
$ws.Range("E3").Value = '  -4.45%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.63%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0981'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.00'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.24'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.58%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.43%  '
$ws.Range("D15").Value = '2.574.59'
$ws.Range("E15").Value = '  -4.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.30%  '
$ws.Range("D18").Value = '2.244.40'
$ws.Range("E18").Value = '  -4.18%  '
$ws.Range("D19").Value = '42.016.25'
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").Value = '0.0₃0974'
$ws.Range("E20").Value = '  -4.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.96%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.49%  '
$ws.Range("E29").Value = '  -3.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.53'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.63%  '
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("E33").Value = '  -6.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0715'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.64%  '
$ws.Range("E38").Value = '  -5.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.68%  '
$ws.Range("E40").Value = '  -6.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0264'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.47'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("B45").Value = 'BitTorrent-New'
$ws.Range("C45").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D45").Value = '0.0₃0164'
$ws.Range("E45").Value = '  +27.71%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.101'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.187'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.14%  '
$ws.Range("B49").Value = 'SynthetixNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.74%  '
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.65%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.28%  '
